$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "CATEGORIES"
$ws.Range("B2").Value = "Customers also Viewed"
$ws.Range("B3").Value = "BEAUTY TO GO"
$ws.Range("B4").Value = "Terms & Conditions"
$ws.Range("B5").Value = "How does the delivery process work?"
$ws.Range("B6").Value = "Privacy Policy"
$ws.Range("B7").Value = "Whey Protein"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "3"
$ws.Range("B9").Value = "Chennai"
$ws.Range("B10").Value = "Sign in"
